$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.222.43'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '2.215.05'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '296.65'
$ws.Range("E5").Value = '  +1.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '88.00'
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("E7").Value = '  +0.43%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.472'
$ws.Range("E9").Value = '  -0.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '52.67'
$ws.Range("E10").Value = '  +7.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '30.95'
$ws.Range("E11").Value = '  +1.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0783'
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.41'
$ws.Range("E14").Value = '  -0.96%  '
$ws.Range("D15").Value = '2.554.64'
$ws.Range("E15").Value = '  -0.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.90'
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("D17").Value = '2.192.69'
$ws.Range("E17").Value = '  -1.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.738'
$ws.Range("E18").Value = '  +1.18%  '
$ws.Range("D19").Value = '40.125.21'
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.38'
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.79'
$ws.Range("E22").Value = '  -0.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.81'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '236.03'
$ws.Range("E24").Value = '  -0.54%  '
$ws.Range("E25").Value = '  +0.09%  '
$ws.Range("E26").Value = '  +0.97%  '
$ws.Range("E27").Value = '  -0.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.33'
$ws.Range("E28").Value = '  +2.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.34'
$ws.Range("E29").Value = '  +1.07%  '
$ws.Range("E30").Value = '  -5.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.27'
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.25'
$ws.Range("E32").Value = '  +1.21%  '
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.99'
$ws.Range("E34").Value = '  +0.55%  '
$ws.Range("E35").Value = '  +3.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0716'
$ws.Range("E36").Value = '  -0.62%  '
$ws.Range("E37").Value = '  -0.74%  '
$ws.Range("E38").Value = '  +1.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.102'
$ws.Range("E39").Value = '  +3.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.74'
$ws.Range("E40").Value = '  +2.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.65'
$ws.Range("E41").Value = '  -1.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.84'
$ws.Range("E42").Value = '  -1.29%  '
$ws.Range("D43").Value = '2.064.28'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.28'
$ws.Range("E44").Value = '  +4.80%  '
$ws.Range("E45").Value = '  +0.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.01'
$ws.Range("E46").Value = '  +1.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.84'
$ws.Range("E47").Value = '  +6.32%  '
$ws.Range("E48").Value = '  -11.16%  '
$ws.Range("D49").Value = '2.428.20'
$ws.Range("E49").Value = '  -0.37%  '
$ws.Range("E50").Value = '  +2.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.47'
$ws.Range("E51").Value = '  +0.27%  '
